$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$ws.Range("B9").Value = "Alvearie Team"

# Old duplicate "Contact" row (A10/B10) becomes the new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate "Contact" row (old row 11); this shifts
# everything below up by one, giving the final A1:B20 range.
$ws.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition columns now describe this specific
# extension instead of the generic Extension placeholder text.
$ws2.Range("K2").Value = "Employee Wage Amount"
$ws2.Range("L2").Value = "Wage amount of the employee for the time period represented by the Wage Basis field"
